$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column D (Modules_required), shifting
# Modules_required/Courses_required/Spots_available/Duration one column
# to the right (D->E, E->F, F->G, G->H).
$ws.Columns.Item(4).Insert()

# New header for the inserted "Location" column
$ws.Cells.Item(1, 4).Value = "Location"

# New dummy "Location" data for rows 2-12
$locations = @(
    "123 Dummy road",
    "124 Dummy road",
    "125 Dummy road",
    "126 Dummy road",
    "127 Dummy road",
    "128 Dummy road",
    "129 Dummy road",
    "130 Dummy road",
    "131 Dummy road",
    "132 Dummy road",
    "133 Dummy road"
)

for ($i = 0; $i -lt $locations.Length; $i++) {
    $ws.Cells.Item($i + 2, 4).Value = $locations[$i]
}

# Resize the URL column (now narrower relative to its content after the
# insert) and the new Location column to fit their text.
$ws.Columns.Item(3).ColumnWidth = 11.09
$ws.Columns.Item(4).ColumnWidth = 7.59

# Match the cell selection left behind in the saved workbook.
$ws.Range("J9").Select()
